# Populate the "Configuration" sheet with the Browser/URL/Button/Label table
# and a hyperlink on the URL cell (Chrome -> https://www.seznam.cz), matching
# the functional-testing configuration added by the commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: A1=Browser, B1=URL, C1=Button, D1=Label
$ws.Range("B1").Value = "URL"
$ws.Range("C1").Value = "Button"
$ws.Range("D1").Value = "Label"
$ws.Range("A1").Value = "Browser"

# Data row: A2=Chrome, B2=https://www.seznam.cz (as a hyperlink)
$ws.Range("A2").Value = "Chrome"
$ws.Range("B2").Value = "https://www.seznam.cz"
$ws.Hyperlinks.Add($ws.Range("B2"), "https://www.seznam.cz")

# Give the whole data row the Hyperlink look (matches A2:D2 styling in the file)
$ws.Range("A2:D2").Style = "Hyperlink"

# Leave the cursor on the hyperlink cell, as in the saved workbook
$ws.Range("B2").Select() | Out-Null
